$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KB")

# Fix the column B (max_pH) boundary values so each band ends exactly where
# the next one begins (e.g. 0.49 -> 0.5, 0.99 -> 1, 2.49 -> 2.5, etc.)
for ($row = 2; $row -le 29; $row++) {
    $ws.Cells.Item($row, 2).Value = 0.5 + ($row - 2) * 0.5
}

# The old row 30 (14 - 14.49 band) is now redundant since row 29 already
# extends up to pH 14, so remove it entirely.
$ws.Rows.Item(30).Delete()

# Update the remembered selection to match the author's final cursor position.
$ws.Range("G30").Select()
